# Barkley Sockeye stock-recruit infilled workbook update
# 1. "metadata" sheet: insert a "comments" row before the "H_cv" row (old row 13),
#    pushing H_cv / S_cv rows down one and leaving a "!definition required!" placeholder.
# 2. "S-R data" sheet: insert a new "comments" column before the H_cv column (old column L),
#    pushing H_cv / S_cv data right one column. Populate the new comments column with
#    source notes for the historical Hucuktlis (HED) data, and backfill a few newly
#    available H / I (run size / harvest) values for 1926-1934.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "metadata"
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("metadata")

# Push the H_cv / S_cv definition rows down by one and insert the new "comments" row.
$wsMeta.Rows(13).Insert()
$wsMeta.Range("A13").Value = "comments"
$wsMeta.Range("B13").Value = "!definition required!"

# ---------------------------------------------------------------------------
# Sheet 2: "S-R data"
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("S-R data")

# Push H_cv (old L) / S_cv (old M) one column to the right and insert the new
# "comments" column in their place. Excel's column insert automatically carries
# the old L values into M and the old M values into N for every row, so the
# only remaining work is: (a) the new header, (b) the comment text for rows
# that need it, and (c) the handful of exceptions noted below.
$wsData.Columns("L").Insert()
$wsData.Range("L1").Value = "comments"

$STARR = "data from Starr, P. J., A. T. Charles and M. A. Henderson 1984. Reconstruction of British Columbia sockeye salmon (Oncorhynchus nerka) stocks: 1970-1982. Can. MS. Rep. Fish. Aquat. Sci"
$HENDERSON = "data from draft Henderson paper"
$HYATT = "data from Hyatt and Steer 1987"

# Historical Hucuktlis (HED) rows 96-99 (1977-1980): Starr et al. 1984 reconstruction.
for ($r = 96; $r -le 99; $r++) {
  $wsData.Cells.Item($r, 12).Value = $STARR
}

# Historical Hucuktlis (HED) rows 100-127 (1981-2008): draft Henderson paper.
for ($r = 100; $r -le 127; $r++) {
  $wsData.Cells.Item($r, 12).Value = $HENDERSON
}

# Historical Hucuktlis (HED) rows 143-194 (1918-1969): Hyatt and Steer 1987.
for ($r = 143; $r -le 194; $r++) {
  $wsData.Cells.Item($r, 12).Value = $HYATT
}

# Historical Hucuktlis (HED) rows 195-201 (1970-1976): Starr et al. 1984 reconstruction.
for ($r = 195; $r -le 201; $r++) {
  $wsData.Cells.Item($r, 12).Value = $STARR
}

# Newly infilled run-size (H) / harvest (I) values for 1926-1933 (rows 151-153, 155-159),
# plus the associated H_cv (now column M) override of 0.05 for those same rows.
$wsData.Cells.Item(151, 8).Value = 100412
$wsData.Cells.Item(151, 9).Value = 35412
$wsData.Cells.Item(151, 13).Value = 0.05

$wsData.Cells.Item(152, 8).Value = 95669
$wsData.Cells.Item(152, 9).Value = 25669
$wsData.Cells.Item(152, 13).Value = 0.05

$wsData.Cells.Item(153, 8).Value = 85000
$wsData.Cells.Item(153, 9).Value = 15000
$wsData.Cells.Item(153, 13).Value = 0.05

$wsData.Cells.Item(155, 8).Value = 50600
$wsData.Cells.Item(155, 9).Value = 10600
$wsData.Cells.Item(155, 13).Value = 0.05

$wsData.Cells.Item(156, 8).Value = 70260
$wsData.Cells.Item(156, 9).Value = 20260
$wsData.Cells.Item(156, 13).Value = 0.05

$wsData.Cells.Item(157, 8).Value = 63000
$wsData.Cells.Item(157, 9).Value = 28000
$wsData.Cells.Item(157, 13).Value = 0.05

$wsData.Cells.Item(158, 8).Value = 28000
$wsData.Cells.Item(158, 9).Value = 20500
$wsData.Cells.Item(158, 13).Value = 0.05

$wsData.Cells.Item(159, 8).Value = 21500
$wsData.Cells.Item(159, 9).Value = 6500
$wsData.Cells.Item(159, 13).Value = 0.05
